# Geofabric-IUCNGET workbook touch-up:
#  - extend the "rdfs:label" (col N) formula on the SSSOM sheet so it also
#    appends the data row number (ROW(B2)-1, ROW(B3)-1, ...)
#  - switch the active tab / selection from "header" to "SSSOM"
#  - leave everything else (data columns A-M, "header" sheet data) untouched

$wb = $excel.ActiveWorkbook

$wsSssom = $wb.Worksheets.Item("SSSOM")

# --- Update the rdfs:label formulas in column N ------------------------
# N2 is a standalone formula; N3:N22 is a shared-formula block (Excel will
# recreate that sharing when a single Formula is assigned across the range).
$wsSssom.Range("N2").Formula = '=CONCAT(B2, " - mapping to IUCN GET - ", ROW(B2)-1 )'
$wsSssom.Range("N3:N22").Formula = '=CONCAT(B3, " - mapping to IUCN GET - ", ROW(B3)-1 )'

# --- Move the active tab / selection from "header" to "SSSOM" ----------
$wsSssom.Activate()
$wsSssom.Range("N3:N22").Select()

# Best-effort: scroll the view so row 7 / column F sit at the top-left of
# the visible pane (matches the author's on-screen scroll position).
try {
    $excel.ActiveWindow.ScrollRow = 7
    $excel.ActiveWindow.ScrollColumn = 6
} catch {
}
